$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RDFModel")

# --- Pass 1: uppercase the Subject column (A) entity identifiers, top to bottom ---
$ws.Range("A2").Value = "PERSON1"
$ws.Range("A3").Value = "PERSON1"
$ws.Range("A4").Value = "PERSON1"
$ws.Range("A5").Value = "PERSON1"

$ws.Range("A6").Value = "STUDY1"
$ws.Range("A7").Value = "STUDY1"

$ws.Range("A8").Value = "TREAT1"
$ws.Range("A9").Value = "TREAT1"

$ws.Range("A10").Value = "PERSON2"
$ws.Range("A11").Value = "PERSON2"
$ws.Range("A12").Value = "PERSON2"
$ws.Range("A13").Value = "PERSON2"

$ws.Range("A14").Value = "PROTOCOL1"

$ws.Range("A15").Value = "STUDY1"
$ws.Range("A16").Value = "STUDY1"
$ws.Range("A17").Value = "STUDY1"

$ws.Range("A18").Value = "TREAT2"

# --- Pass 2: fix up the matching Object (C) references in the second block ---
$ws.Range("C10").Value = "STUDY1"
$ws.Range("C13").Value = "TREAT2"
$ws.Range("C16").Value = "PROTOCOL1"
$ws.Range("C17").Value = "TREAT2"

# Move the selection to B25
$ws.Range("B25").Select()

# Protect the sheet (matches sheetProtection element added to sheet1.xml)
$ws.Protect("password")
